$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 52
$ws_ALC.Range("H52").Value = 47621724
$ws_ALC.Range("I52").Value = 3049.5
$ws_ALC.Range("J52").Value = 52634216
$ws_ALC.Range("K52").Value = 9148.5
$ws_ALC.Range("L52").Value = 157902648
$ws_ALC.Range("M52").Value = -8988.5
$ws_ALC.Range("N52").Value = -157902968

# ALC row 64
$ws_ALC.Range("H64").Value = 2621.0962
$ws_ALC.Range("I64").Value = 2617.261
$ws_ALC.Range("J64").Value = 2624.138
$ws_ALC.Range("K64").Value = 2617.261
$ws_ALC.Range("L64").Value = 2624.138
$ws_ALC.Range("M64").Value = -2369.261
$ws_ALC.Range("N64").Value = -3120.138

# ALC row 67
$ws_ALC.Range("H67").Value = 2621.0962
$ws_ALC.Range("I67").Value = 2617.261
$ws_ALC.Range("J67").Value = 2624.138
$ws_ALC.Range("K67").Value = 2617.261
$ws_ALC.Range("L67").Value = 2624.138
$ws_ALC.Range("M67").Value = -1759.261
$ws_ALC.Range("N67").Value = -4340.138

# ALC row 132
$ws_ALC.Range("H132").Value = 5621538.5
$ws_ALC.Range("I132").Value = 8067925.5
$ws_ALC.Range("J132").Value = 3909.4814
$ws_ALC.Range("K132").Value = 24203776.5
$ws_ALC.Range("L132").Value = 11728.4442
$ws_ALC.Range("M132").Value = -24201246.5
$ws_ALC.Range("N132").Value = -16788.4442

# ALC row 137
$ws_ALC.Range("H137").Value = 1077.9286
$ws_ALC.Range("I137").Value = 952.17645
$ws_ALC.Range("J137").Value = 1612.375
$ws_ALC.Range("K137").Value = 2856.52935
$ws_ALC.Range("L137").Value = 4837.125
$ws_ALC.Range("M137").Value = -306.5293500000002
$ws_ALC.Range("N137").Value = -9937.125

# ALC row 138
$ws_ALC.Range("H138").Value = 1914.386
$ws_ALC.Range("I138").Value = 1084.4062
$ws_ALC.Range("J138").Value = 2976.76
$ws_ALC.Range("K138").Value = 3253.2186
$ws_ALC.Range("L138").Value = 8930.280000000001
$ws_ALC.Range("M138").Value = 1886.7814
$ws_ALC.Range("N138").Value = -19210.28

# ALC row 141
$ws_ALC.Range("H141").Value = 1493.091
$ws_ALC.Range("I141").Value = 912.11475
$ws_ALC.Range("K141").Value = 2736.34425
$ws_ALC.Range("M141").Value = 2443.65575

# ARM row 37
$ws_ARM.Range("H37").Value = 10108.667
$ws_ARM.Range("J37").Value = 10530.4
$ws_ARM.Range("L37").Value = 10530.4
$ws_ARM.Range("N37").Value = -11076.4

# ARM row 61
$ws_ARM.Range("H61").Value = 1287.3864
$ws_ARM.Range("I61").Value = 892.7143
$ws_ARM.Range("K61").Value = 892.7143
$ws_ARM.Range("M61").Value = -680.7143

# ARM row 136
$ws_ARM.Range("H136").Value = 1287.3864
$ws_ARM.Range("I136").Value = 892.7143
$ws_ARM.Range("K136").Value = 2678.1429
$ws_ARM.Range("M136").Value = -128.1428999999998

# BSM row 134
$ws_BSM.Range("H134").Value = 2756.3157
$ws_BSM.Range("I134").Value = 738.9796
$ws_BSM.Range("K134").Value = 2216.9388
$ws_BSM.Range("M134").Value = 318.0612000000001

# CRP row 31
$ws_CRP.Range("H31").Value = 1406.3137
$ws_CRP.Range("I31").Value = 1113.1875
$ws_CRP.Range("J31").Value = 1540.3143
$ws_CRP.Range("K31").Value = 1113.1875
$ws_CRP.Range("L31").Value = 1540.3143
$ws_CRP.Range("M31").Value = -818.1875
$ws_CRP.Range("N31").Value = -2130.3143

# CRP row 34
$ws_CRP.Range("H34").Value = 1406.3137
$ws_CRP.Range("I34").Value = 1113.1875
$ws_CRP.Range("J34").Value = 1540.3143
$ws_CRP.Range("K34").Value = 1113.1875
$ws_CRP.Range("L34").Value = 1540.3143
$ws_CRP.Range("M34").Value = -911.1875
$ws_CRP.Range("N34").Value = -1944.3143

# CRP row 58
$ws_CRP.Range("H58").Value = 13889676
$ws_CRP.Range("I58").Value = 17857890
$ws_CRP.Range("J58").Value = 925.875
$ws_CRP.Range("K58").Value = 17857890
$ws_CRP.Range("L58").Value = 925.875
$ws_CRP.Range("M58").Value = -17857687
$ws_CRP.Range("N58").Value = -1331.875

# CRP row 62
$ws_CRP.Range("H62").Value = 3639.8572
$ws_CRP.Range("I62").Value = 2367.5715
$ws_CRP.Range("J62").Value = 4912.143
$ws_CRP.Range("K62").Value = 2367.5715
$ws_CRP.Range("L62").Value = 4912.143
$ws_CRP.Range("M62").Value = -1743.5715
$ws_CRP.Range("N62").Value = -6160.143

# CRP row 65
$ws_CRP.Range("H65").Value = 3639.8572
$ws_CRP.Range("I65").Value = 2367.5715
$ws_CRP.Range("J65").Value = 4912.143
$ws_CRP.Range("K65").Value = 11837.8575
$ws_CRP.Range("L65").Value = 24560.715
$ws_CRP.Range("M65").Value = -8717.8575
$ws_CRP.Range("N65").Value = -30800.715

# CRP row 132
$ws_CRP.Range("H132").Value = 7577153
$ws_CRP.Range("I132").Value = 1042.4138
$ws_CRP.Range("J132").Value = 22224300
$ws_CRP.Range("K132").Value = 3127.2414
$ws_CRP.Range("L132").Value = 66672900
$ws_CRP.Range("M132").Value = -597.2413999999999
$ws_CRP.Range("N132").Value = -66677960

# CRP row 134
$ws_CRP.Range("H134").Value = 14706758
$ws_CRP.Range("I134").Value = 863
$ws_CRP.Range("J134").Value = 83334264
$ws_CRP.Range("K134").Value = 2589
$ws_CRP.Range("L134").Value = 250002792
$ws_CRP.Range("M134").Value = -54
$ws_CRP.Range("N134").Value = -250007862

# CRP row 136
$ws_CRP.Range("H136").Value = 13889676
$ws_CRP.Range("I136").Value = 17857890
$ws_CRP.Range("J136").Value = 925.875
$ws_CRP.Range("K136").Value = 53573670
$ws_CRP.Range("L136").Value = 2777.625
$ws_CRP.Range("M136").Value = -53571120
$ws_CRP.Range("N136").Value = -7877.625

# CUL row 34
$ws_CUL.Range("H34").Value = 575.1429000000001
$ws_CUL.Range("I34").Value = 563
$ws_CUL.Range("J34").Value = 580
$ws_CUL.Range("K34").Value = 1689
$ws_CUL.Range("L34").Value = 1740
$ws_CUL.Range("M34").Value = -1605
$ws_CUL.Range("N34").Value = -1908

# CUL row 55
$ws_CUL.Range("H55").Value = 912.5
$ws_CUL.Range("I55").Value = 150
$ws_CUL.Range("J55").Value = 1166.6666
$ws_CUL.Range("K55").Value = 450
$ws_CUL.Range("L55").Value = 3499.9998
$ws_CUL.Range("M55").Value = -273
$ws_CUL.Range("N55").Value = -3853.9998

# CUL row 76
$ws_CUL.Range("H76").Value = 1250
$ws_CUL.Range("I76").Value = 500
$ws_CUL.Range("K76").Value = 1500
$ws_CUL.Range("M76").Value = -1117

# CUL row 79
$ws_CUL.Range("H79").Value = 1250
$ws_CUL.Range("I79").Value = 500
$ws_CUL.Range("K79").Value = 1500
$ws_CUL.Range("M79").Value = -174

# GSM row 132
$ws_GSM.Range("H132").Value = 7397.5713
$ws_GSM.Range("I132").Value = 4473.1924
$ws_GSM.Range("J132").Value = 15845.777
$ws_GSM.Range("K132").Value = 13419.5772
$ws_GSM.Range("L132").Value = 47537.331
$ws_GSM.Range("M132").Value = -10889.5772
$ws_GSM.Range("N132").Value = -52597.331

# LTW row 136
$ws_LTW.Range("H136").Value = 28389974
$ws_LTW.Range("I136").Value = 10353689
$ws_LTW.Range("J136").Value = 166668160
$ws_LTW.Range("K136").Value = 31061067
$ws_LTW.Range("L136").Value = 500004480
$ws_LTW.Range("M136").Value = -31058517
$ws_LTW.Range("N136").Value = -500009580

# WVR row 132
$ws_WVR.Range("H132").Value = 3869.9827
$ws_WVR.Range("I132").Value = 2249.0571
$ws_WVR.Range("J132").Value = 6336.609
$ws_WVR.Range("K132").Value = 6747.1713
$ws_WVR.Range("L132").Value = 19009.827
$ws_WVR.Range("M132").Value = -4217.1713
$ws_WVR.Range("N132").Value = -24069.827

# WVR row 136
$ws_WVR.Range("H136").Value = 1507.3077
$ws_WVR.Range("I136").Value = 825.8333
$ws_WVR.Range("J136").Value = 3778.889
$ws_WVR.Range("K136").Value = 2477.4999
$ws_WVR.Range("L136").Value = 11336.667
$ws_WVR.Range("M136").Value = 72.5001000000002
$ws_WVR.Range("N136").Value = -16436.667
